# "do dl buoi hoc" - update seed data paths to include subfolder prefixes
# (baigiang/ for lecture PDFs, anh/ for image paths) and repoint the
# bai_giang rows at khoa_hoc_id 7, then leave the workbook focused on the
# bai_giang sheet the way the author last left it.

$wb = $excel.ActiveWorkbook

# --- bai_giang sheet: filebaigiang now lives under baigiang/, and every
#     lecture row now belongs to khoa_hoc_id 7 instead of 1 ---
$wsBaiGiang = $wb.Worksheets.Item("bai_giang")
$wsBaiGiang.Range("C2").Value = "baigiang/bai_giang_1_1.pdf"
$wsBaiGiang.Range("D2").Value = 7
$wsBaiGiang.Range("C3").Value = "baigiang/bai_giang_1_2.pdf"
$wsBaiGiang.Range("D3").Value = 7
$wsBaiGiang.Range("C4").Value = "baigiang/bai_giang_1_3.pdf"
$wsBaiGiang.Range("D4").Value = 7
$wsBaiGiang.Range("C5").Value = "baigiang/bai_giang_1_4.pdf"
$wsBaiGiang.Range("D5").Value = 7
$wsBaiGiang.Range("C6").Value = "baigiang/bai_giang_1_5.pdf"
$wsBaiGiang.Range("D6").Value = 7

# --- hoc_cu sheet: equipment images now live under anh/ ---
$wsHocCu = $wb.Worksheets.Item("hoc_cu")
$wsHocCu.Range("B2").Value = "anh/hoccu/1.png"
$wsHocCu.Range("B3").Value = "anh/hoccu/2.png"
$wsHocCu.Range("B4").Value = "anh/hoccu/3.png"

# --- hinh_anh_khoa_hoc sheet: course images now live under anh/ ---
$wsHinhAnh = $wb.Worksheets.Item("hinh_anh_khoa_hoc")
$wsHinhAnh.Range("A2").Value = "anh/khoahoc/mncn/1.png"
$wsHinhAnh.Range("A3").Value = "anh/khoahoc/mncn/2.png"
$wsHinhAnh.Range("A4").Value = "anh/khoahoc/mncn/3.png"
$wsHinhAnh.Range("A5").Value = "anh/khoahoc/belamgame/1.png"
$wsHinhAnh.Range("A6").Value = "anh/khoahoc/belamgame/2.png"
$wsHinhAnh.Range("A7").Value = "anh/khoahoc/belamgame/3.png"

# --- column width tweaks left behind by the author's last edit ---
$wsHinhAnh.Columns.Item(1).ColumnWidth = 35.8
$wsBaiGiang.Columns.Item(2).ColumnWidth = 23.65
$wsBaiGiang.Columns.Item(3).ColumnWidth = 21.3

# --- restore/update each sheet's last-used selection ---
$wsKhoaHoc = $wb.Worksheets.Item("khoa_hoc")
$wsKhoaHoc.Range("D36").Select() | Out-Null

$wsHinhAnh.Range("B16").Select() | Out-Null

$wsHocCu.Range("B9").Select() | Out-Null

$wsBaiTap = $wb.Worksheets.Item("bai_tap")
$wsBaiTap.Range("E15").Select() | Out-Null

# bai_giang is the sheet the author left active/selected
$wsBaiGiang.Range("D12").Select() | Out-Null
